$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / email-list cells - direct value assignment is safe (no
# numeric/date/percent auto-conversion risk).
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System"
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G7").Value = "Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("H16").Value = "111/251"
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"

# Percentage-looking text cells ("25.5%") - a plain Value assignment would
# be auto-interpreted by Excel as a numeric percent (changing the cell's
# type/style to a number format). Route the literal text through a
# formula first, then collapse the formula back down to a plain value via
# copy / paste-special so the stored cell keeps its original text type and
# style (reading .Value back into a variable is not reliable here, so we
# stay entirely inside the COM object model for the round-trip).
$ws.Range("L10").Formula = '="25.5%"'
$ws.Range("L10").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4163) | Out-Null

$ws.Range("S15").Formula = '="25.5%"'
$ws.Range("S15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = $false
